# Auto-generated edit script: update "want to go" counts (column F) for
# several events, plus mark one event as sold out (G21 on sheet 1).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 3379
$ws.Range("F3").Value = 797
$ws.Range("F4").Value = 2442
$ws.Range("F5").Value = 523
$ws.Range("F6").Value = 442
$ws.Range("F7").Value = 266
$ws.Range("F9").Value = 408
$ws.Range("F10").Value = 1141
$ws.Range("F11").Value = 497
$ws.Range("F12").Value = 259
$ws.Range("F13").Value = 103
$ws.Range("F14").Value = 312
$ws.Range("F15").Value = 5166
$ws.Range("F16").Value = 36
$ws.Range("F17").Value = 1438
$ws.Range("F18").Value = 3787
$ws.Range("F19").Value = 364
$ws.Range("F21").Value = 310
$ws.Range("G21").Value = "已售罄"
$ws.Range("F22").Value = 4186
$ws.Range("F23").Value = 5560
$ws.Range("F25").Value = 998
$ws.Range("F26").Value = 601
$ws.Range("F27").Value = 3513
$ws.Range("F28").Value = 423
$ws.Range("F29").Value = 58
$ws.Range("F30").Value = 159
$ws.Range("F31").Value = 103
$ws.Range("F32").Value = 934
$ws.Range("F33").Value = 1268
$ws.Range("F34").Value = 87
$ws.Range("F35").Value = 119
$ws.Range("F36").Value = 1502
$ws.Range("F37").Value = 167
$ws.Range("F38").Value = 1520
$ws.Range("F39").Value = 87
$ws.Range("F40").Value = 982
$ws.Range("F41").Value = 1013
$ws.Range("F42").Value = 564
$ws.Range("F43").Value = 71
$ws.Range("F44").Value = 104
$ws.Range("F45").Value = 2601
$ws.Range("F46").Value = 104
$ws.Range("F49").Value = 3804

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 1095
$ws.Range("F9").Value = 12
$ws.Range("F16").Value = 2
$ws.Range("F22").Value = 55

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 3337

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 3337
$ws.Range("F3").Value = 797
$ws.Range("F4").Value = 2442
$ws.Range("F5").Value = 523
$ws.Range("F6").Value = 442
$ws.Range("F7").Value = 266
$ws.Range("F8").Value = 1095
$ws.Range("F10").Value = 408
$ws.Range("F11").Value = 1141
$ws.Range("F12").Value = 497
$ws.Range("F13").Value = 259
$ws.Range("F14").Value = 103
$ws.Range("F15").Value = 312
$ws.Range("F16").Value = 5166
$ws.Range("F18").Value = 1438
$ws.Range("F19").Value = 4186
$ws.Range("F20").Value = 5560
$ws.Range("F22").Value = 998
$ws.Range("F23").Value = 601
$ws.Range("F24").Value = 3513
$ws.Range("F25").Value = 423
$ws.Range("F26").Value = 58
$ws.Range("F27").Value = 159
$ws.Range("F28").Value = 103
$ws.Range("F29").Value = 934
$ws.Range("F30").Value = 1268
$ws.Range("F31").Value = 87
$ws.Range("F32").Value = 119
$ws.Range("F33").Value = 1502
$ws.Range("F34").Value = 167
$ws.Range("F35").Value = 1520
$ws.Range("F36").Value = 2
$ws.Range("F37").Value = 982
$ws.Range("F39").Value = 564
$ws.Range("F41").Value = 71
$ws.Range("F42").Value = 55
$ws.Range("F43").Value = 2601
$ws.Range("F45").Value = 104
$ws.Range("F49").Value = 3804
